$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": row 2 holds the handoff/handback datetimes for the
# 7f9414b6-... zh-cn.xlf file. Update Correspond Handoff/Handback Datetime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 10:46:55"
$wsZhCn.Range("H2").Value = "2016-03-24 10:47:36"

# Sheet "de-de": row 2 holds the handoff/handback datetimes for the
# 7f9414b6-... de-de.xlf file. Update Correspond Handoff/Handback Datetime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 10:47:04"
$wsDeDe.Range("H2").Value = "2016-03-24 10:47:50"
